$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("H6").Value = 2.88
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("S6").Value = 1.67
$ws.Range("T6").Value = 2.1
$ws.Range("AH6").Value = 12
$ws.Range("AJ6").Value = 34
$ws.Range("AK6").Value = 34
$ws.Range("AL6").Value = 51
$ws.Range("AT6").Value = 2.1
$ws.Range("BA6").Value = 126
$ws.Range("H12").Value = 3.7
$ws.Range("J12").Value = 2.3
$ws.Range("K12").Value = 2.22
$ws.Range("O12").Value = 1.21
$ws.Range("P12").Value = 3.55
$ws.Range("R12").Value = 2.02
$ws.Range("V12").Value = 2.05
$ws.Range("Y12").Value = 8
$ws.Range("AD12").Value = 7.3
$ws.Range("AG12").Value = 14
$ws.Range("AP12").Value = 16
$ws.Range("AQ12").Value = 28
$ws.Range("AT12").Value = 3
$ws.Range("AW12").Value = 5.9
$ws.Range("AX12").Value = 22
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 2.63
$ws.Range("AA13").Value = 17
$ws.Range("AB13").Value = 21
$ws.Range("I14").Value = 2.9
$ws.Range("K14").Value = 1.83
$ws.Range("L14").Value = 4
$ws.Range("M14").Value = 1.13
$ws.Range("N14").Value = 6
$ws.Range("S14").Value = 1.62
$ws.Range("T14").Value = 2.2
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 1.57
$ws.Range("X14").Value = 10
$ws.Range("Z14").Value = 23
$ws.Range("AC14").Value = 6
$ws.Range("AE14").Value = 21
$ws.Range("AG14").Value = 6.5
$ws.Range("AK14").Value = 34
$ws.Range("AP14").Value = 34
$ws.Range("AS14").Value = 351
$ws.Range("AT14").Value = 2.2
$ws.Range("AU14").Value = 9.5
$ws.Range("Q23").Value = 2.15
$ws.Range("R23").Value = 1.67
$ws.Range("O25").Value = 1.25
$ws.Range("P25").Value = 3.75
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93
$ws.Range("H27").Value = 3.15
$ws.Range("I27").Value = 3.05
$ws.Range("J27").Value = 2.82
$ws.Range("L27").Value = 3.6
$ws.Range("O27").Value = 1.34
$ws.Range("P27").Value = 3
$ws.Range("S27").Value = 1.42
$ws.Range("T27").Value = 2.67
$ws.Range("W27").Value = 7.4
$ws.Range("AA27").Value = 18.5
$ws.Range("AF27").Value = 65
$ws.Range("AG27").Value = 9
$ws.Range("AH27").Value = 16
$ws.Range("AK27").Value = 27
$ws.Range("AL27").Value = 35
$ws.Range("AO27").Value = 11.75
$ws.Range("AP27").Value = 20
$ws.Range("AR27").Value = 80
$ws.Range("AT27").Value = 2.67
$ws.Range("AV27").Value = 60
$ws.Range("AY27").Value = 24
$ws.Range("O30").Value = 1.16
$ws.Range("P30").Value = 4.65
$ws.Range("V30").Value = 2.5
$ws.Range("X30").Value = 11.75
$ws.Range("BA30").Value = 75
$ws.Range("H31").Value = 3.65
$ws.Range("I31").Value = 4.25
$ws.Range("K31").Value = 2.2
$ws.Range("Q31").Value = 1.82
$ws.Range("U31").Value = 1.78
$ws.Range("V31").Value = 1.93
$ws.Range("W31").Value = 7.2
$ws.Range("AB31").Value = 25
$ws.Range("AG31").Value = 12.5
$ws.Range("AH31").Value = 25
$ws.Range("AJ31").Value = 70
$ws.Range("AL31").Value = 45
$ws.Range("AN31").Value = 3.6
$ws.Range("AO31").Value = 8.25
$ws.Range("AP31").Value = 17.5
$ws.Range("G34").Value = 2.45
$ws.Range("I34").Value = 2.65
$ws.Range("K34").Value = 2.12
$ws.Range("L34").Value = 3.2
$ws.Range("T34").Value = 2.75
$ws.Range("W34").Value = 7.9
$ws.Range("X34").Value = 12
$ws.Range("Z34").Value = 26
$ws.Range("AB34").Value = 32
$ws.Range("AD34").Value = 6.3
$ws.Range("AG34").Value = 8.5
$ws.Range("AH34").Value = 13.5
$ws.Range("AK34").Value = 22
$ws.Range("AN34").Value = 4.4
$ws.Range("AT34").Value = 2.75
$ws.Range("AU34").Value = 7
$ws.Range("AW34").Value = 4.65
$ws.Range("AX34").Value = 14
$ws.Range("AY34").Value = 21
$ws.Range("AZ34").Value = 60
$ws.Range("BA34").Value = 90
$ws.Range("BB34").Value = 250
$ws.Range("G35").Value = 1.8
$ws.Range("I35").Value = 4.45
$ws.Range("J35").Value = 2.4
$ws.Range("K35").Value = 2.05
$ws.Range("L35").Value = 4.7
$ws.Range("P35").Value = 3
$ws.Range("Q35").Value = 1.9
$ws.Range("R35").Value = 1.8
$ws.Range("T35").Value = 2.52
$ws.Range("V35").Value = 1.88
$ws.Range("W35").Value = 6.8
$ws.Range("X35").Value = 8.5
$ws.Range("AB35").Value = 25
$ws.Range("AC35").Value = 9.25
$ws.Range("AE35").Value = 14
$ws.Range("AF35").Value = 65
$ws.Range("AG35").Value = 12.5
$ws.Range("AH35").Value = 27
$ws.Range("AK35").Value = 45
$ws.Range("AL35").Value = 45
$ws.Range("AM35").Value = 500
$ws.Range("AO35").Value = 9
$ws.Range("AP35").Value = 18
$ws.Range("AQ35").Value = 32
$ws.Range("AR35").Value = 65
$ws.Range("AT35").Value = 2.47
$ws.Range("AU35").Value = 7
$ws.Range("AV35").Value = 65
$ws.Range("AY35").Value = 29
$ws.Range("BA35").Value = 175
$ws.Range("G36").Value = 2.9
$ws.Range("I36").Value = 2.5
$ws.Range("J36").Value = 3.45
$ws.Range("T36").Value = 2.45
$ws.Range("AG36").Value = 8.5
$ws.Range("AY36").Value = 20
$ws.Range("AZ36").Value = 60
